$d = $word.ActiveDocument

$replacements = @(
    @("767×5=", "738×5="),
    @("959×5=", "518×6="),
    @("601×9=", "461×3="),
    @("336×6=", "501×2="),
    @("596×3=", "179×4="),
    @("979×8=", "405×4="),
    @("583×2=", "615×3="),
    @("130×4=", "266×7="),
    @("527×7=", "922×7="),
    @("781×7=", "522×6="),
    @("646×2=", "527×6="),
    @("455×4=", "171×5="),
    @("466×9=", "269×3="),
    @("978×2=", "140×4="),
    @("984×7=", "340×2="),
    @("279×9=", "965×8="),
    @("466×7=", "233×2="),
    @("867×8=", "950×9="),
    @("220×2=", "904×3="),
    @("316×7=", "445×5="),
    @("613×7=", "285×2="),
    @("319×3=", "976×9="),
    @("811×3=", "184×8="),
    @("617×9=", "926×9="),
    @("267×2=", "652×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
